# "update for BAR chart"
#
# Root cause of the whole diff: the sample reading in `_input!B3` changes
# from 47952 to 200. `_input!D12` (=$B3), `_input!I12` (=D12-E12) and
# `_input!J12` (=D12-F12) / `_input!K12` (=SUMIF ...) all derive from it, and
# `trend!C28/D28/H28/I28/M28` pull those cached results in turn - Excel's
# recalculation engine keeps all of that in sync once we poke the root cell.
#
# On top of that, `_input!E12` / `_input!F12` used to be hard literal 0s;
# they become real formulas pointing at `_input2!$B3` / `_input3!$B3` (the
# "yesterday" / "last week" helper sheets), and those helper sheets' sample
# value is reset from 45583 / 136 down to 0. Selecting B3 on each helper
# sheet mirrors the <selection> bookkeeping Excel leaves behind when a user
# last clicks on that cell.

$wb = $excel.ActiveWorkbook

# --- _input2 ("yesterday") sample value -> 0 -----------------------------
$wsInput2 = $wb.Worksheets.Item("_input2")
$wsInput2.Range("B3").Value = 0
$wsInput2.Range("B3").Select()

# --- _input3 ("last week") sample value -> 0 ------------------------------
$wsInput3 = $wb.Worksheets.Item("_input3")
$wsInput3.Range("B3").Value = 0
$wsInput3.Range("B3").Select()

# --- _input: wire E12/F12 to the helper sheets, then update the root value
$wsInput = $wb.Worksheets.Item("_input")
$wsInput.Range("E12").Formula = '=_input2!$B3'
$wsInput.Range("F12").Formula = '=_input3!$B3'
$wsInput.Range("B3").Value = 200

# Recalculate everything so cached <v> results (trend!C28/D28/H28/I28/M28,
# _input!D12/I12/J12/K12, etc.) are refreshed before save.
$excel.CalculateFullRebuild()

# Restore the original active sheet/selection on "trend" (it was the
# selected tab before this edit).
$wsTrend = $wb.Worksheets.Item("trend")
$wsTrend.Activate()
